# Updated symbol list on Wed Feb 15 19:44:00 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for the
# coin rows that moved since the last scrape. Values are written with a
# leading apostrophe so Excel keeps them as plain text (matching the
# workbook's existing inline-string / text cells) instead of silently
# re-interpreting numeric- or percent-looking strings as real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'305.11"
$ws.Range("E2").Value  = "'2.55%"

$ws.Range("D3").Value  = "'44.06"
$ws.Range("E3").Value  = "'6.29%"

$ws.Range("D4").Value  = "'5.124"
$ws.Range("E4").Value  = "'2.33%"

$ws.Range("D5").Value  = "'0.07833"
$ws.Range("E5").Value  = "'3.90%"

$ws.Range("D6").Value  = "'4.436"
$ws.Range("E6").Value  = "'1.66%"

$ws.Range("D7").Value  = "'1.615"
$ws.Range("E7").Value  = "'2.79%"

$ws.Range("D8").Value  = "'1.056"
$ws.Range("E8").Value  = "'13.48%"

$ws.Range("D9").Value  = "'0.1306"
$ws.Range("E9").Value  = "'7.80%"

$ws.Range("D10").Value = "'0.1872"
$ws.Range("E10").Value = "'2.36%"

$ws.Range("D11").Value = "'0.09164"
$ws.Range("E11").Value = "'3.43%"

$ws.Range("E12").Value = "'1.69%"

$ws.Range("D13").Value = "'0.1047"
$ws.Range("E13").Value = "'-0.70%"

$ws.Range("D14").Value = "'0.001297"
$ws.Range("E14").Value = "'1.16%"

$ws.Range("D15").Value = "'0.005865"
$ws.Range("E15").Value = "'-0.17%"

$ws.Range("E17").Value = "'0.77%"

$ws.Range("E18").Value = "'-2.40%"

$ws.Range("D19").Value = "'0.3357"
$ws.Range("E19").Value = "'0.69%"

$ws.Range("D20").Value = "'8.040"
$ws.Range("E20").Value = "'1.05%"

$ws.Range("D21").Value = "'0.1379"
$ws.Range("E21").Value = "'-2.64%"

$ws.Range("E22").Value = "'-5.16%"

$ws.Range("D23").Value = "'0.04187"
$ws.Range("E23").Value = "'3.30%"

$ws.Range("D24").Value = "'0.001273"
$ws.Range("E24").Value = "'0.70%"

$ws.Range("D25").Value = "'0.004486"
$ws.Range("E25").Value = "'14.61%"

$ws.Range("D26").Value = "'0.0001341"
$ws.Range("E26").Value = "'9.15%"

$ws.Range("D38").Value = "'0.02574"
$ws.Range("E38").Value = "'6.38%"

$ws.Range("D39").Value = "'0.05347"
$ws.Range("E39").Value = "'2.49%"

$ws.Range("D40").Value = "'0.005595"
$ws.Range("E40").Value = "'-4.89%"

$ws.Range("D41").Value = "'0.007742"
$ws.Range("E41").Value = "'-0.55%"

$ws.Range("D42").Value = "'0.1379"
$ws.Range("E42").Value = "'3.57%"

$ws.Range("D43").Value = "'0.007332"
$ws.Range("E43").Value = "'-0.27%"

$ws.Range("D44").Value = "'0.008328"

$ws.Range("D45").Value = "'0.3020"
$ws.Range("E45").Value = "'1.36%"

$ws.Range("D46").Value = "'0.00006682"
$ws.Range("E46").Value = "'6.11%"

$ws.Range("E47").Value = "'-0.52%"

$ws.Range("D48").Value = "'0.04413"
$ws.Range("E48").Value = "'-1.36%"

$ws.Range("D49").Value = "'0.003971"
$ws.Range("E49").Value = "'-5.33%"

$ws.Range("D50").Value = "'0.00002086"
$ws.Range("E50").Value = "'-0.52%"

$ws.Range("D51").Value = "'0.0001987"
$ws.Range("E51").Value = "'-0.52%"
